# Refresh the cryptos price list (Price + Volume(1h) columns) and shift the
# tail of the coin ranking (RenderToken dropped out of the table, every
# following coin moved up one slot, Maker entered at the bottom).
# Values are written with a leading apostrophe so Excel stores them as the
# same literal text (e.g. "1.00", "0.0229") the source data used, instead of
# re-parsing them as numbers and dropping meaningful trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.994.10"
$ws.Range("E2").Value = "'  -1.11%  "
$ws.Range("D3").Value = "'2.632.98"
$ws.Range("E3").Value = "'  +1.27%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "'  +0.18%  "
$ws.Range("D5").Value = "'514.32"
$ws.Range("E5").Value = "'  +0.28%  "
$ws.Range("D6").Value = "'144.47"
$ws.Range("E6").Value = "'  -1.09%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "'  -0.16%  "
$ws.Range("E8").Value = "'  +1.56%  "
$ws.Range("D9").Value = "'2.660.23"
$ws.Range("E9").Value = "'  +2.25%  "
$ws.Range("D10").Value = "'6.30"
$ws.Range("E10").Value = "'  +1.00%  "
$ws.Range("E11").Value = "'  +2.49%  "
$ws.Range("E12").Value = "'  +0.47%  "
$ws.Range("E13").Value = "'  -1.34%  "
$ws.Range("D14").Value = "'3.123.49"
$ws.Range("E14").Value = "'  +2.05%  "
$ws.Range("D15").Value = "'59.023.58"
$ws.Range("E15").Value = "'  -1.09%  "
$ws.Range("D16").Value = "'21.11"
$ws.Range("E16").Value = "'  +1.09%  "
$ws.Range("E17").Value = "'  +1.16%  "
$ws.Range("D18").Value = "'2.659.89"
$ws.Range("E18").Value = "'  +1.92%  "
$ws.Range("E19").Value = "'  -0.40%  "
$ws.Range("D20").Value = "'343.96"
$ws.Range("E20").Value = "'  +1.46%  "
$ws.Range("D21").Value = "'10.40"
$ws.Range("E21").Value = "'  +1.67%  "
$ws.Range("E22").Value = "'  +1.38%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("E23").Value = "'  +0.30%  "
$ws.Range("D24").Value = "'61.01"
$ws.Range("E24").Value = "'  +0.62%  "
$ws.Range("E25").Value = "'  +2.18%  "
$ws.Range("D26").Value = "'2.753.50"
$ws.Range("E26").Value = "'  +0.85%  "
$ws.Range("D27").Value = "'0.994"
$ws.Range("E27").Value = "'  -0.67%  "
$ws.Range("D28").Value = "'0.161"
$ws.Range("E28").Value = "'  +2.45%  "
$ws.Range("E29").Value = "'  +2.41%  "
$ws.Range("E30").Value = "'  +3.15%  "
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "'  -0.10%  "
$ws.Range("E32").Value = "'  +9.37%  "
$ws.Range("E33").Value = "'  +1.20%  "
$ws.Range("D34").Value = "'18.92"
$ws.Range("E34").Value = "'  +1.38%  "
$ws.Range("D35").Value = "'149.50"
$ws.Range("E35").Value = "'  -0.52%  "
$ws.Range("E36").Value = "'  +12.73%  "
$ws.Range("E37").Value = "'  +4.44%  "
$ws.Range("D38").Value = "'1.15"
$ws.Range("E38").Value = "'  +3.18%  "
$ws.Range("E39").Value = "'  +2.25%  "
$ws.Range("D40").Value = "'36.47"
$ws.Range("E40").Value = "'  -0.19%  "
$ws.Range("E41").Value = "'  +3.86%  "
$ws.Range("E42").Value = "'  +0.56%  "
$ws.Range("D43").Value = "'280.49"
$ws.Range("E43").Value = "'  -0.91%  "
$ws.Range("D44").Value = "'0.616"
$ws.Range("E44").Value = "'  -0.54%  "
$ws.Range("D45").Value = "'0.994"
$ws.Range("E45").Value = "'  -0.57%  "
$ws.Range("D46").Value = "'0.0984"
$ws.Range("E46").Value = "'  -0.36%  "
$ws.Range("E47").Value = "'  +3.05%  "
$ws.Range("E48").Value = "'  -0.44%  "
$ws.Range("B49").Value = "'WhiteBITCoin"
$ws.Range("C49").Value = "'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").Value = "'10.27"
$ws.Range("E49").Value = "'  -1.02%  "
$ws.Range("B50").Value = "'VeChain"
$ws.Range("C50").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0229"
$ws.Range("E50").Value = "'  +0.24%  "
$ws.Range("B51").Value = "'Maker"
$ws.Range("C51").Value = "'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").Value = "'1.982.08"
$ws.Range("E51").Value = "'  +2.65%  "
